$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: new comment text (index 5 string)
$ws.Range("B2").Value = "Just schema consistent. `n6/6 old/new PAs`nadjascent locations.`nSmall board"

# Add new rows 3 and 4
$ws.Range("A3").Value = "new2"
$ws.Range("A4").Value = "new3"

$ws.Range("B3").Value = "Larger board. But same as 48652"
$ws.Range("B4").Value = "Larger board. 8/6 adjascent!"

# Adjust column B width and selection
$ws.Columns.Item(2).ColumnWidth = 29.0
$ws.Range("B15").Select()
